$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 674.34784
$ws.Range("I6").Value = 313.52942
$ws.Range("J6").Value = 1696.6666
$ws.Range("K6").Value = 940.58826
$ws.Range("L6").Value = 5089.9998
$ws.Range("M6").Value = -828.58826
$ws.Range("N6").Value = -5313.9998

$ws.Range("H18").Value = 1293.75
$ws.Range("I18").Value = 309.25
$ws.Range("K18").Value = 309.25
$ws.Range("M18").Value = -25.25

$ws.Range("H38").Value = 511
$ws.Range("I38").Value = 28.714285
$ws.Range("J38").Value = 1636.3334
$ws.Range("K38").Value = 86.142855
$ws.Range("L38").Value = 4909.0002
$ws.Range("M38").Value = 285.857145
$ws.Range("N38").Value = -5653.0002

$ws.Range("H51").Value = 4715.8335
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 4715.8335
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 4715.8335
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -5683.8335

$ws.Range("H52").Value = 600
$ws.Range("I52").Value = 600
$ws.Range("K52").Value = 1800
$ws.Range("M52").Value = -1640

$ws.Range("H58").Value = 4012.6924
$ws.Range("I58").Value = 144.16667
$ws.Range("J58").Value = 7328.5713
$ws.Range("K58").Value = 432.50001
$ws.Range("L58").Value = 21985.7139
$ws.Range("M58").Value = -282.50001
$ws.Range("N58").Value = -22285.7139

$ws.Range("H62").Value = 7435720.5
$ws.Range("J62").Value = 35362.375
$ws.Range("L62").Value = 35362.375
$ws.Range("N62").Value = -36610.375

$ws.Range("H65").Value = 7435720.5
$ws.Range("J65").Value = 35362.375
$ws.Range("L65").Value = 176811.875
$ws.Range("N65").Value = -183051.875

$ws.Range("H68").Value = 22000
$ws.Range("J68").Value = 22000
$ws.Range("L68").Value = 22000
$ws.Range("N68").Value = -23498

$ws.Range("H71").Value = 22000
$ws.Range("J71").Value = 22000
$ws.Range("L71").Value = 66000
$ws.Range("N71").Value = -73488

$ws.Range("H86").Value = 3042.5454
$ws.Range("I86").Value = 2534
$ws.Range("J86").Value = 3466.3333
$ws.Range("K86").Value = 2534
$ws.Range("L86").Value = 3466.3333
$ws.Range("M86").Value = -1411
$ws.Range("N86").Value = -5712.3333

$ws.Range("H88").Value = 2472.2942
$ws.Range("I88").Value = 1377.25
$ws.Range("J88").Value = 2809.2307
$ws.Range("K88").Value = 1377.25
$ws.Range("L88").Value = 2809.2307
$ws.Range("M88").Value = -971.25
$ws.Range("N88").Value = -3621.2307

$ws.Range("H89").Value = 3042.5454
$ws.Range("I89").Value = 2534
$ws.Range("J89").Value = 3466.3333
$ws.Range("K89").Value = 12670
$ws.Range("L89").Value = 17331.6665
$ws.Range("M89").Value = -7054
$ws.Range("N89").Value = -28563.6665

$ws.Range("H91").Value = 2472.2942
$ws.Range("I91").Value = 1377.25
$ws.Range("J91").Value = 2809.2307
$ws.Range("K91").Value = 1377.25
$ws.Range("L91").Value = 2809.2307
$ws.Range("M91").Value = 26.75
$ws.Range("N91").Value = -5617.2307

$ws.Range("H98").Value = 512869.7
$ws.Range("I98").Value = 593037.75
$ws.Range("J98").Value = 5138.6665
$ws.Range("K98").Value = 593037.75
$ws.Range("L98").Value = 5138.6665
$ws.Range("M98").Value = -591539.75
$ws.Range("N98").Value = -8134.6665

$ws.Range("H106").Value = 37038304
$ws.Range("I106").Value = 37038304
$ws.Range("K106").Value = 37038304
$ws.Range("M106").Value = -37037673

$ws.Range("H122").Value = 512869.7
$ws.Range("I122").Value = 593037.75
$ws.Range("J122").Value = 5138.6665
$ws.Range("K122").Value = 1779113.25
$ws.Range("L122").Value = 15415.9995
$ws.Range("M122").Value = -1776663.25
$ws.Range("N122").Value = -20315.9995

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws.Range("H128").Value = 98888
$ws.Range("J128").Value = 98888
$ws.Range("L128").Value = 98888
$ws.Range("N128").Value = -108848

$ws.Range("H130").Value = 49000
$ws.Range("J130").Value = 49000
$ws.Range("L130").Value = 49000
$ws.Range("N130").Value = -59040

$ws.Range("H137").Value = 26317120
$ws.Range("I137").Value = 32259050
$ws.Range("K137").Value = 96777150
$ws.Range("M137").Value = -96774600

$ws.Range("H138").Value = 8098144.5
$ws.Range("J138").Value = 10641116
$ws.Range("L138").Value = 31923348
$ws.Range("N138").Value = -31933628

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3228.2104
$ws.Range("I61").Value = 2652.3845
$ws.Range("J61").Value = 4475.8335
$ws.Range("K61").Value = 2652.3845
$ws.Range("L61").Value = 4475.8335
$ws.Range("M61").Value = -2440.3845
$ws.Range("N61").Value = -4899.8335

$ws.Range("H136").Value = 3228.2104
$ws.Range("I136").Value = 2652.3845
$ws.Range("J136").Value = 4475.8335
$ws.Range("K136").Value = 7957.1535
$ws.Range("L136").Value = 13427.5005
$ws.Range("M136").Value = -5407.1535
$ws.Range("N136").Value = -18527.5005

$ws.Range("H139").Value = 51206.668
$ws.Range("J139").Value = 51206.668
$ws.Range("L139").Value = 51206.668
$ws.Range("N139").Value = -61486.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 525004.5
$ws.Range("I14").Value = 50000
$ws.Range("J14").Value = 1000009
$ws.Range("K14").Value = 50000
$ws.Range("L14").Value = 1000009
$ws.Range("M14").Value = -49828
$ws.Range("N14").Value = -1000353

$ws.Range("H107").Value = 2925
$ws.Range("I107").Value = 2903
$ws.Range("J107").Value = 3013
$ws.Range("K107").Value = 2903
$ws.Range("L107").Value = 3013
$ws.Range("M107").Value = -983
$ws.Range("N107").Value = -6853

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 1000000
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()

$ws.Range("H31").Value = 3946.2163
$ws.Range("I31").Value = 1734.0869
$ws.Range("J31").Value = 7580.4287
$ws.Range("K31").Value = 1734.0869
$ws.Range("L31").Value = 7580.4287
$ws.Range("M31").Value = -1439.0869
$ws.Range("N31").Value = -8170.4287

$ws.Range("H34").Value = 3946.2163
$ws.Range("I34").Value = 1734.0869
$ws.Range("J34").Value = 7580.4287
$ws.Range("K34").Value = 1734.0869
$ws.Range("L34").Value = 7580.4287
$ws.Range("M34").Value = -1532.0869
$ws.Range("N34").Value = -7984.4287

$ws.Range("H105").Value = 958.0909
$ws.Range("I105").Value = 961
$ws.Range("K105").Value = 961
$ws.Range("M105").Value = 786

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 10102845
$ws.Range("I131").Value = 380
$ws.Range("J131").Value = 11906857
$ws.Range("K131").Value = 1140
$ws.Range("L131").Value = 35720571
$ws.Range("M131").Value = 3900
$ws.Range("N131").Value = -35730651

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2219.8572
$ws.Range("I82").Value = 2135.8
$ws.Range("K82").Value = 2135.8
$ws.Range("M82").Value = -1774.8

$ws.Range("H85").Value = 2219.8572
$ws.Range("I85").Value = 2135.8
$ws.Range("K85").Value = 2135.8
$ws.Range("M85").Value = -887.8000000000002

$ws.Range("H95").Value = 29500
$ws.Range("J95").Value = 29500
$ws.Range("L95").Value = 29500
$ws.Range("N95").Value = -34992

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 5820
$ws.Range("I51").Value = 3500
$ws.Range("J51").Value = 6400
$ws.Range("K51").Value = 3500
$ws.Range("L51").Value = 6400
$ws.Range("M51").Value = -2990
$ws.Range("N51").Value = -7420

$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
